$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1751
$ws.Range("F5").Value = 1320
$ws.Range("F6").Value = 387
$ws.Range("F9").Value = 739
$ws.Range("F11").Value = 543
$ws.Range("F15").Value = 3061
$ws.Range("F16").Value = 2678
$ws.Range("F18").Value = 34
$ws.Range("F20").Value = 327
$ws.Range("F23").Value = 5476
$ws.Range("F25").Value = 1004
$ws.Range("F26").Value = 38
$ws.Range("F27").Value = 69
$ws.Range("F28").Value = 389
$ws.Range("F29").Value = 1158
$ws.Range("F32").Value = 307

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 43
$ws.Range("F11").Value = 1
$ws.Range("F25").Value = 4008

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F9").Value = 1389

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 1751
$ws.Range("F8").Value = 1389
$ws.Range("F12").Value = 1320
$ws.Range("F13").Value = 387
$ws.Range("F15").Value = 739
$ws.Range("F18").Value = 543
$ws.Range("F20").Value = 3061
$ws.Range("F21").Value = 2678
$ws.Range("F22").Value = 34
$ws.Range("F24").Value = 327
$ws.Range("F25").Value = 43
$ws.Range("F28").Value = 5476
$ws.Range("F30").Value = 1004
$ws.Range("F32").Value = 38
$ws.Range("F33").Value = 69
$ws.Range("F34").Value = 389
$ws.Range("F41").Value = 1158
$ws.Range("F49").Value = 307
